$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 43, pushing the existing rows 43-51
# down to 45-53 (dimension grows from A1:R51 to A1:R53).
$ws.Rows.Item(43).Insert()
$ws.Rows.Item(43).Insert()

# New row 43: Alcachofa, Española, Primera, from Provincia de Limarí
$ws.Cells.Item(43, 1).Value = 5
$ws.Cells.Item(43, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(43, 3).Value = "Maule"
$ws.Cells.Item(43, 4).Value = 44466
$ws.Cells.Item(43, 5).Value = 7
$ws.Cells.Item(43, 6).Value = 100112013
$ws.Cells.Item(43, 7).Value = "Alcachofa"
$ws.Cells.Item(43, 8).Value = "Española"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 300
$ws.Cells.Item(43, 11).Value = 11000
$ws.Cells.Item(43, 12).Value = 11000
$ws.Cells.Item(43, 13).Value = 11000
$ws.Cells.Item(43, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(43, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(43, 16).Value = 367
$ws.Cells.Item(43, 17).Value = 30
$ws.Cells.Item(43, 18).Value = "Hortaliza"

# New row 44: Alcachofa, Madrigal, Primera, from Provincia de Limarí
$ws.Cells.Item(44, 1).Value = 5
$ws.Cells.Item(44, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(44, 3).Value = "Maule"
$ws.Cells.Item(44, 4).Value = 44466
$ws.Cells.Item(44, 5).Value = 7
$ws.Cells.Item(44, 6).Value = 100112013
$ws.Cells.Item(44, 7).Value = "Alcachofa"
$ws.Cells.Item(44, 8).Value = "Madrigal"
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 300
$ws.Cells.Item(44, 11).Value = 10000
$ws.Cells.Item(44, 12).Value = 10000
$ws.Cells.Item(44, 13).Value = 10000
$ws.Cells.Item(44, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(44, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(44, 16).Value = 250
$ws.Cells.Item(44, 17).Value = 40
$ws.Cells.Item(44, 18).Value = "Hortaliza"
